$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update instructor names to append " ( بكالوريوس )" qualifier
$ws.Range("C2").Value = "م. رائد محمد مراد ( بكالوريوس )"
$ws.Range("C16").Value = "أ. إسماعيل دواس ( بكالوريوس )"
$ws.Range("C20").Value = "د. مراد أبومنسي ( بكالوريوس )"
$ws.Range("C23").Value = "د. محمد السردي ( بكالوريوس )"
$ws.Range("C24").Value = "أ. عبد الله أبو قاسم  ( بكالوريوس )"

# Adjust column widths (target stored widths: 17.5 and 28.875 character units;
# the host quantizes ColumnWidth to 1/7-character steps, so feed the nearest
# value that rounds to the closest achievable stored width)
$ws.Columns.Item(2).ColumnWidth = 16.857142857142858
$ws.Columns.Item(3).ColumnWidth = 28.142857142857142

# Update selection
$ws.Range("C24").Select()
